$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new text would otherwise be auto-coerced to a number by Excel;
# set them as Text first, assign the literal string, then restore the default style
# so the cell formatting matches the original (unstyled) cells.
$priceCells = @("D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D18","D19","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D5").Value = "1.000"
$ws.Range("D6").Value = "289.85"
$ws.Range("D7").Value = "0.3681"
$ws.Range("D8").Value = "49.49"
$ws.Range("D9").Value = "0.3393"
$ws.Range("D10").Value = "1.172"
$ws.Range("D11").Value = "0.07627"
$ws.Range("D13").Value = "21.44"
$ws.Range("D14").Value = "6.072"
$ws.Range("D15").Value = "6.943"
$ws.Range("D16").Value = "0.00001140"
$ws.Range("D18").Value = "89.36"
$ws.Range("D19").Value = "0.06755"
$ws.Range("D21").Value = "6.258"
$ws.Range("D22").Value = "16.61"
$ws.Range("D23").Value = "0.5309"
$ws.Range("D24").Value = "12.02"
$ws.Range("D26").Value = "2.385"
$ws.Range("D27").Value = "2.975"
$ws.Range("D28").Value = "20.03"
$ws.Range("D29").Value = "145.98"
$ws.Range("D30").Value = "4.982"
$ws.Range("D31").Value = "126.00"
$ws.Range("D33").Value = "1.049"
$ws.Range("D34").Value = "6.311"
$ws.Range("D35").Value = "2.005"
$ws.Range("D36").Value = "10.34"
$ws.Range("D37").Value = "0.08455"
$ws.Range("D38").Value = "0.02542"
$ws.Range("D39").Value = "0.2335"
$ws.Range("D40").Value = "0.06580"
$ws.Range("D41").Value = "5.575"
$ws.Range("D42").Value = "11.80"
$ws.Range("D44").Value = "0.6398"
$ws.Range("D45").Value = "14.33"
$ws.Range("D46").Value = "1.000"
$ws.Range("D47").Value = "0.6011"
$ws.Range("D48").Value = "3.746"
$ws.Range("D50").Value = "1.262"
$ws.Range("D51").Value = "123.51"
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

# Price cells whose new text already looks non-numeric (contains two or more dots),
# so Excel keeps them as plain text without any extra handling.
$ws.Range("D2").Value = "22.404.43"
$ws.Range("D3").Value = "1.575.09"
$ws.Range("D17").Value = "1.569.25"
$ws.Range("D25").Value = "22.405.75"
$ws.Range("D32").Value = "1.743.16"

# Volume(1h) percentage cells; the surrounding spaces keep Excel from treating them as numbers.
$ws.Range("E2").Value = "  -3.95%  "
$ws.Range("E3").Value = "  -3.29%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("E6").Value = "  -2.55%  "
$ws.Range("E7").Value = "  -2.01%  "
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("E9").Value = "  -3.11%  "
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("E11").Value = "  -4.61%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("E15").Value = "  -3.64%  "
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("E17").Value = "  -3.79%  "
$ws.Range("E18").Value = "  -5.82%  "
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  -5.59%  "
$ws.Range("E22").Value = "  -3.63%  "
$ws.Range("E23").Value = "  -6.86%  "
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("E25").Value = "  -3.99%  "
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("E29").Value = "  -3.70%  "
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("E31").Value = "  -4.00%  "
$ws.Range("E32").Value = "  -3.83%  "
$ws.Range("E33").Value = "  +9.17%  "
$ws.Range("E34").Value = "  -6.45%  "
$ws.Range("E35").Value = "  -5.19%  "
$ws.Range("E36").Value = "  -6.98%  "
$ws.Range("E37").Value = "  -2.66%  "
$ws.Range("E38").Value = "  -4.96%  "
$ws.Range("E39").Value = "  -2.97%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("E42").Value = "  -7.31%  "
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("E44").Value = "  -5.61%  "
$ws.Range("E45").Value = "  -6.56%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("E47").Value = "  -4.18%  "
$ws.Range("E48").Value = "  -3.55%  "
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("E50").Value = "  +4.82%  "
$ws.Range("E51").Value = "  -2.12%  "
